$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("errors")

# The window moved slightly on screen between edits (xWindow/yWindow in the
# saved workbookView). Reflect that via the window position.
$win = $wb.Windows.Item(1)
$win.Left = 7300
$win.Top = 1620

# Populate row 6 with the new "objectIsLocked" error entry (account-system
# "object is locked" error code 505).
# Write the Chinese message first so it lands at the lower shared-string index,
# matching the order new strings were appended to the shared strings table.
$ws.Range("C6").Value = "对象被锁定"
$ws.Range("A6").Value = "objectIsLocked"
$ws.Range("B6").Value = 505

# Move the active selection to A7, matching the post-edit cursor position
$ws.Range("A7").Select()
